$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 750
$ws.Range("I18").Value = 750
$ws.Range("K18").Value = 750
$ws.Range("M18").Value = -466

$ws.Range("H28").Value = 2258.6365
$ws.Range("I28").Value = 2258.6365
$ws.Range("K28").Value = 2258.6365
$ws.Range("M28").Value = -1773.6365

$ws.Range("H58").Value = 281.06668
$ws.Range("I58").Value = 170.6923
$ws.Range("J58").Value = 998.5
$ws.Range("K58").Value = 512.0769
$ws.Range("L58").Value = 2995.5
$ws.Range("M58").Value = -362.0769
$ws.Range("N58").Value = -3295.5

$ws.Range("H62").Value = 2393
$ws.Range("I62").Value = 2192.125
$ws.Range("K62").Value = 2192.125
$ws.Range("M62").Value = -1568.125

$ws.Range("H65").Value = 2393
$ws.Range("I65").Value = 2192.125
$ws.Range("K65").Value = 10960.625
$ws.Range("M65").Value = -7840.625

$ws.Range("H70").Value = 5855.8335
$ws.Range("I70").Value = 5911.5
$ws.Range("J70").Value = 5800.1665
$ws.Range("K70").Value = 17734.5
$ws.Range("L70").Value = 17400.4995
$ws.Range("M70").Value = -17464.5
$ws.Range("N70").Value = -17940.4995

$ws.Range("H73").Value = 5855.8335
$ws.Range("I73").Value = 5911.5
$ws.Range("J73").Value = 5800.1665
$ws.Range("K73").Value = 17734.5
$ws.Range("L73").Value = 17400.4995
$ws.Range("M73").Value = -16798.5
$ws.Range("N73").Value = -19272.4995

$ws.Range("H88").Value = 4611.923
$ws.Range("J88").Value = 5282.364
$ws.Range("L88").Value = 5282.364
$ws.Range("N88").Value = -6094.364

$ws.Range("H91").Value = 4611.923
$ws.Range("J91").Value = 5282.364
$ws.Range("L91").Value = 5282.364
$ws.Range("N91").Value = -8090.364

$ws.Range("H113").Value = 5949.25
$ws.Range("I113").Value = 1399.5
$ws.Range("K113").Value = 1399.5
$ws.Range("M113").Value = 1854.5

$ws.Range("H129").Value = 1986.5625
$ws.Range("J129").Value = 2500
$ws.Range("L129").Value = 7500
$ws.Range("N129").Value = -17500

$ws.Range("H132").Value = 4493.2856
$ws.Range("I132").Value = 4659.7085
$ws.Range("J132").Value = 3494.75
$ws.Range("K132").Value = 13979.1255
$ws.Range("L132").Value = 10484.25
$ws.Range("M132").Value = -11449.1255
$ws.Range("N132").Value = -15544.25

$ws.Range("H136").Value = 90780
$ws.Range("J136").Value = 90780
$ws.Range("L136").Value = 90780
$ws.Range("N136").Value = -100980

$ws.Range("H137").Value = 3723.9092
$ws.Range("I137").Value = 4242.5
$ws.Range("K137").Value = 12727.5
$ws.Range("M137").Value = -10177.5

$ws.Range("H138").Value = 2696.276
$ws.Range("J138").Value = 3431.75
$ws.Range("L138").Value = 10295.25
$ws.Range("N138").Value = -20575.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4556.25
$ws.Range("I32").Value = 4556.25
$ws.Range("K32").Value = 4556.25
$ws.Range("M32").Value = -4269.25

$ws.Range("H56").Value = 26000
$ws.Range("I56").Value = 26000
$ws.Range("K56").Value = 26000
$ws.Range("M56").Value = -25258

$ws.Range("H131").Value = 79999
$ws.Range("J131").Value = 79999
$ws.Range("L131").Value = 79999
$ws.Range("N131").Value = -90079

$ws.Range("H132").Value = 3625.2559
$ws.Range("I132").Value = 3628.238
$ws.Range("K132").Value = 10884.714
$ws.Range("M132").Value = -8354.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 275
$ws.Range("I22").Value = 275
$ws.Range("K22").Value = 275
$ws.Range("M22").Value = -102

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H134").Value = 2576.7368
$ws.Range("I134").Value = 1984.1765
$ws.Range("K134").Value = 5952.529500000001
$ws.Range("M134").Value = -3417.529500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 15384810
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H24").Value = 15384810
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H41").Value = 12535.875
$ws.Range("J41").Value = 28250
$ws.Range("L41").Value = 28250
$ws.Range("N41").Value = -29106

$ws.Range("H58").Value = 3343
$ws.Range("I58").Value = 3343
$ws.Range("K58").Value = 3343
$ws.Range("M58").Value = -3140

$ws.Range("H129").Value = 94000
$ws.Range("J129").Value = 94000
$ws.Range("L129").Value = 94000
$ws.Range("N129").Value = -104000

$ws.Range("H132").Value = 3254.2222
$ws.Range("I132").Value = 3261
$ws.Range("K132").Value = 9783
$ws.Range("M132").Value = -7253

$ws.Range("H134").Value = 7520887
$ws.Range("I134").Value = 7938436.5
$ws.Range("K134").Value = 23815309.5
$ws.Range("M134").Value = -23812774.5

$ws.Range("H136").Value = 3343
$ws.Range("I136").Value = 3343
$ws.Range("K136").Value = 10029
$ws.Range("M136").Value = -7479

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 116.4375
$ws.Range("I2").Value = 78.111115
$ws.Range("J2").Value = 165.71428
$ws.Range("K2").Value = 468.66669
$ws.Range("L2").Value = 994.28568
$ws.Range("M2").Value = -355.66669
$ws.Range("N2").Value = -1220.28568

$ws.Range("H12").Value = 135.11111
$ws.Range("I12").Value = 239.8
$ws.Range("J12").Value = 4.25
$ws.Range("K12").Value = 719.4000000000001
$ws.Range("L12").Value = 12.75
$ws.Range("M12").Value = -546.4000000000001
$ws.Range("N12").Value = -358.75

$ws.Range("H13").Value = 21.833334
$ws.Range("J13").Value = 15
$ws.Range("L13").Value = 45
$ws.Range("N13").Value = -381

$ws.Range("H38").Value = 28
$ws.Range("I38").Value = 36.666668
$ws.Range("J38").Value = 17.6
$ws.Range("K38").Value = 110.000004
$ws.Range("L38").Value = 52.8
$ws.Range("M38").Value = 236.999996
$ws.Range("N38").Value = -746.8

$ws.Range("H70").Value = 13398.5
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 13398.5
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H116").Value = 1842.6666
$ws.Range("I116").Value = 1842.6666
$ws.Range("K116").Value = 5527.9998
$ws.Range("M116").Value = -2085.9998

$ws.Range("H131").Value = 2007.9474
$ws.Range("I131").Value = 1861.9
$ws.Range("J131").Value = 2170.2222
$ws.Range("K131").Value = 5585.700000000001
$ws.Range("L131").Value = 6510.6666
$ws.Range("M131").Value = -545.7000000000007
$ws.Range("N131").Value = -16590.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 24388.777
$ws.Range("I46").Value = 15583.333
$ws.Range("J46").Value = 41999.668
$ws.Range("K46").Value = 15583.333
$ws.Range("L46").Value = 41999.668
$ws.Range("M46").Value = -15427.333
$ws.Range("N46").Value = -42311.668

$ws.Range("H57").Value = 25000
$ws.Range("I57").Value = 25000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 25000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -24180
$ws.Range("N57").ClearContents()

$ws.Range("H135").Value = 39424
$ws.Range("I135").Value = 39424
$ws.Range("K135").Value = 39424
$ws.Range("M135").Value = -34354

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 993.2857
$ws.Range("I16").Value = 992
$ws.Range("J16").Value = 1001
$ws.Range("K16").Value = 992
$ws.Range("L16").Value = 1001
$ws.Range("M16").Value = -822
$ws.Range("N16").Value = -1341

$ws.Range("H22").Value = 3933.3076
$ws.Range("J22").Value = 6861.5
$ws.Range("L22").Value = 6861.5
$ws.Range("N22").Value = -7451.5

$ws.Range("H27").Value = 3933.3076
$ws.Range("J27").Value = 6861.5
$ws.Range("L27").Value = 6861.5
$ws.Range("N27").Value = -7075.5

$ws.Range("H61").Value = 3396.2307
$ws.Range("I61").Value = 3218.3635
$ws.Range("K61").Value = 3218.3635
$ws.Range("M61").Value = -3016.3635

$ws.Range("H100").Value = 3575
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H113").Value = 3396.2307
$ws.Range("I113").Value = 3218.3635
$ws.Range("K113").Value = 3218.3635
$ws.Range("M113").Value = -1048.3635

$ws.Range("H136").Value = 1899.5
$ws.Range("I136").Value = 1899.5
$ws.Range("K136").Value = 5698.5
$ws.Range("M136").Value = -3148.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3581.7097
$ws.Range("I132").Value = 2424.64
$ws.Range("J132").Value = 8402.833000000001
$ws.Range("K132").Value = 7273.92
$ws.Range("L132").Value = 25208.499
$ws.Range("M132").Value = -4743.92
$ws.Range("N132").Value = -30268.499

$ws.Range("H136").Value = 5229.5654
$ws.Range("I136").Value = 5674.15
$ws.Range("J136").Value = 2265.6667
$ws.Range("K136").Value = 17022.45
$ws.Range("L136").Value = 6797.000100000001
$ws.Range("M136").Value = -14472.45
$ws.Range("N136").Value = -11897.0001
